$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = "'22.402.39"
$ws.Range('E2').Formula = "'  -0.13%  "
$ws.Range('D3').Formula = "'1.566.86"
$ws.Range('E3').Formula = "'  +0.02%  "
$ws.Range('D4').Formula = "'1.000"
$ws.Range('E4').Formula = "'  -0.16%  "
$ws.Range('E5').Formula = "'  -0.05%  "
$ws.Range('D6').Formula = "'286.70"
$ws.Range('E6').Formula = "'  +0.74%  "
$ws.Range('D7').Formula = "'0.3730"
$ws.Range('E7').Formula = "'  +2.83%  "
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Formula = "'0.3284"
$ws.Range('E8').Formula = "'  -0.94%  "
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Formula = "'45.68"
$ws.Range('E9').Formula = "'  -5.35%  "
$ws.Range('D10').Formula = "'1.144"
$ws.Range('E10').Formula = "'  +2.11%  "
$ws.Range('D11').Formula = "'0.07413"
$ws.Range('E11').Formula = "'  +0.42%  "
$ws.Range('E12').Formula = "'  -0.15%  "
$ws.Range('D13').Formula = "'20.48"
$ws.Range('E13').Formula = "'  -1.06%  "
$ws.Range('D14').Formula = "'5.852"
$ws.Range('E14').Formula = "'  -1.57%  "
$ws.Range('D15').Formula = "'6.846"
$ws.Range('E15').Formula = "'  -0.60%  "
$ws.Range('D16').Formula = "'1.568.08"
$ws.Range('E16').Formula = "'  +0.13%  "
$ws.Range('D17').Formula = "'0.00001099"
$ws.Range('E17').Formula = "'  -0.16%  "
$ws.Range('D18').Formula = "'0.06689"
$ws.Range('E18').Formula = "'  -0.33%  "
$ws.Range('D19').Formula = "'85.89"
$ws.Range('E19').Formula = "'  -1.49%  "
$ws.Range('D20').Formula = "'1.000"
$ws.Range('E20').Formula = "'  -0.16%  "
$ws.Range('D21').Formula = "'6.371"
$ws.Range('E21').Formula = "'  +0.14%  "
$ws.Range('D22').Formula = "'16.30"
$ws.Range('E22').Formula = "'  +0.75%  "
$ws.Range('D23').Formula = "'11.77"
$ws.Range('E23').Formula = "'  -1.83%  "
$ws.Range('D24').Formula = "'22.378.97"
$ws.Range('E24').Formula = "'  -0.17%  "
$ws.Range('D25').Formula = "'2.318"
$ws.Range('E25').Formula = "'  -2.06%  "
$ws.Range('D26').Formula = "'2.576"
$ws.Range('E26').Formula = "'  +1.56%  "
$ws.Range('D27').Formula = "'151.88"
$ws.Range('E27').Formula = "'  +0.89%  "
$ws.Range('D28').Formula = "'19.32"
$ws.Range('E28').Formula = "'  -0.25%  "
$ws.Range('D29').Formula = "'4.921"
$ws.Range('E29').Formula = "'  -1.59%  "
$ws.Range('D30').Formula = "'123.70"
$ws.Range('E30').Formula = "'  -0.22%  "
$ws.Range('D31').Formula = "'1.744.50"
$ws.Range('E31').Formula = "'  -0.41%  "
$ws.Range('D32').Formula = "'1.075"
$ws.Range('E32').Formula = "'  +4.15%  "
$ws.Range('D33').Formula = "'1.954"
$ws.Range('E33').Formula = "'  -2.25%  "
$ws.Range('D34').Formula = "'5.954"
$ws.Range('E34').Formula = "'  -1.93%  "
$ws.Range('D35').Formula = "'9.714"
$ws.Range('E35').Formula = "'  +0.64%  "
$ws.Range('D36').Formula = "'0.08270"
$ws.Range('E36').Formula = "'  +0.63%  "
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Formula = "'0.02382"
$ws.Range('E37').Formula = "'  -0.79%  "
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Formula = "'1.301"
$ws.Range('E38').Formula = "'  +1.05%  "
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Formula = "'0.2191"
$ws.Range('E39').Formula = "'  -1.52%  "
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Formula = "'0.06299"
$ws.Range('E40').Formula = "'  -1.57%  "
$ws.Range('D41').Formula = "'5.294"
$ws.Range('E41').Formula = "'  -0.96%  "
$ws.Range('D42').Formula = "'11.13"
$ws.Range('E42').Formula = "'  -0.18%  "
$ws.Range('D43').Formula = "'0.6108"
$ws.Range('E43').Formula = "'  -1.73%  "
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Formula = "'13.80"
$ws.Range('E44').Formula = "'  +0.72%  "
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Formula = "'0.5951"
$ws.Range('E45').Formula = "'  -0.69%  "
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Formula = "'3.747"
$ws.Range('E46').Formula = "'  +0.10%  "
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Formula = "'2.007"
$ws.Range('E47').Formula = "'  -0.79%  "
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Formula = "'124.19"
$ws.Range('E48').Formula = "'  +0.78%  "
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Formula = "'1.176"
$ws.Range('E49').Formula = "'  -2.40%  "
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Formula = "'0.07152"
$ws.Range('E50').Formula = "'  -0.57%  "
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Formula = "'76.17"
$ws.Range('E51').Formula = "'  +1.06%  "
